# Auto-generated PowerShell/COM script to apply F-column ("想去人数") value updates
# across the four worksheets, matching the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7945
$ws.Range("F5").Value = 2154
$ws.Range("F6").Value = 8683
$ws.Range("F8").Value = 97
$ws.Range("F9").Value = 5813
$ws.Range("F10").Value = 67
$ws.Range("F11").Value = 2821
$ws.Range("F12").Value = 1223
$ws.Range("F13").Value = 422
$ws.Range("F16").Value = 641
$ws.Range("F17").Value = 123
$ws.Range("F18").Value = 4040
$ws.Range("F24").Value = 37
$ws.Range("F25").Value = 5883
$ws.Range("F26").Value = 215
$ws.Range("F27").Value = 82
$ws.Range("F28").Value = 292
$ws.Range("F31").Value = 424
$ws.Range("F32").Value = 4296
$ws.Range("F36").Value = 5741
$ws.Range("F37").Value = 85
$ws.Range("F38").Value = 16
$ws.Range("F40").Value = 36
$ws.Range("F41").Value = 3740
$ws.Range("F42").Value = 31
$ws.Range("F43").Value = 45
$ws.Range("F45").Value = 2367
$ws.Range("F49").Value = 42

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 166
$ws.Range("F5").Value = 80
$ws.Range("F10").Value = 136

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1382

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1382
$ws.Range("F4").Value = 7945
$ws.Range("F5").Value = 2154
$ws.Range("F6").Value = 8683
$ws.Range("F7").Value = 97
$ws.Range("F8").Value = 5813
$ws.Range("F9").Value = 67
$ws.Range("F10").Value = 2821
$ws.Range("F11").Value = 1223
$ws.Range("F12").Value = 422
$ws.Range("F15").Value = 166
$ws.Range("F16").Value = 641
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 4040
$ws.Range("F25").Value = 37
$ws.Range("F26").Value = 5883
$ws.Range("F27").Value = 215
$ws.Range("F28").Value = 82
$ws.Range("F31").Value = 424
$ws.Range("F32").Value = 80
$ws.Range("F33").Value = 4296
$ws.Range("F38").Value = 5741
$ws.Range("F39").Value = 85
$ws.Range("F40").Value = 16
$ws.Range("F42").Value = 3740
$ws.Range("F44").Value = 2367
$ws.Range("F47").Value = 42
$ws.Range("F48").Value = 136
